$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D (Price) cells from Excel auto-numeric-coercion so
# values like "0.550" or "38.10" keep their exact text (incl. trailing zeros),
# and multi-dot values parse as text anyway. Apply a Text format first,
# then clear the style back to Normal afterwards so no stray "s" attribute
# is left on the cell (matching the original inlineStr cells with no style).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '42.794.80'
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = '2.556.44'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '304.45'
$ws.Range("E5").Value = '  +1.83%  '
$ws.Range("D6").Value = '98.19'
$ws.Range("E6").Value = '  +4.22%  '
$ws.Range("D7").Value = '0.575'
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '0.550'
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("D10").Value = '37.38'
$ws.Range("E10").Value = '  +3.81%  '
$ws.Range("D11").Value = '0.0808'
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("D12").Value = '0.116'
$ws.Range("E12").Value = '  +7.62%  '
$ws.Range("D13").Value = '7.74'
$ws.Range("E13").Value = '  -0.15%  '
$ws.Range("D14").Value = '2.507.07'
$ws.Range("E14").Value = '  -1.16%  '
$ws.Range("D15").Value = '15.14'
$ws.Range("E15").Value = '  +6.45%  '
$ws.Range("D16").Value = '0.885'
$ws.Range("E16").Value = '  +1.64%  '
$ws.Range("D17").Value = '42.836.87'
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("D18").Value = '13.45'
$ws.Range("E18").Value = '  +7.39%  '
$ws.Range("D19").Value = '0.0₃0988'
$ws.Range("E19").Value = '  +0.74%  '
$ws.Range("D20").Value = '6.66'
$ws.Range("E20").Value = '  -0.18%  '
$ws.Range("D21").Value = '71.66'
$ws.Range("E21").Value = '  -0.49%  '
$ws.Range("D22").Value = '254.73'
$ws.Range("E22").Value = '  -2.67%  '
$ws.Range("E23").Value = '  +1.80%  '
$ws.Range("E24").Value = '  -2.71%  '
$ws.Range("D25").Value = '27.81'
$ws.Range("E25").Value = '  -6.24%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").Value = '10.14'
$ws.Range("E27").Value = '  +0.28%  '
$ws.Range("D28").Value = '38.10'
$ws.Range("E28").Value = '  +2.95%  '
$ws.Range("E29").Value = '  -1.33%  '
$ws.Range("E30").Value = '  +0.54%  '
$ws.Range("D31").Value = '157.54'
$ws.Range("E31").Value = '  +1.95%  '
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("E33").Value = '  +0.66%  '
$ws.Range("D34").Value = '0.0808'
$ws.Range("E34").Value = '  +1.50%  '
$ws.Range("E35").Value = '  -3.02%  '
$ws.Range("D36").Value = '26.50'
$ws.Range("E36").Value = '  +10.25%  '
$ws.Range("D37").Value = '18.66'
$ws.Range("E37").Value = '  +12.06%  '
$ws.Range("E38").Value = '  -0.72%  '
$ws.Range("E39").Value = '  -0.27%  '
$ws.Range("D40").Value = '2.16'
$ws.Range("E40").Value = '  +36.28%  '
$ws.Range("D41").Value = '3.45'
$ws.Range("E41").Value = '  -1.06%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").Value = '0.0305'
$ws.Range("E43").Value = '  -2.57%  '
$ws.Range("D44").Value = '2.076.09'
$ws.Range("E44").Value = '  -0.59%  '
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("D46").Value = '87.37'
$ws.Range("E46").Value = '  +1.69%  '
$ws.Range("D47").Value = '9.13'
$ws.Range("E47").Value = '  +4.82%  '
$ws.Range("D48").Value = '2.804.01'
$ws.Range("E48").Value = '  +0.15%  '
$ws.Range("D49").Value = '74.96'
$ws.Range("E49").Value = '  +7.86%  '
$ws.Range("D50").Value = '103.77'
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("D51").Value = '0.191'
$ws.Range("E51").Value = '  +1.84%  '

$priceRange.Style = "Normal"

